$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3088
$ws1.Range("F5").Value = 1687
$ws1.Range("F6").Value = 2102
$ws1.Range("F12").Value = 443
$ws1.Range("F13").Value = 1147
$ws1.Range("F17").Value = 7509
$ws1.Range("F18").Value = 317
$ws1.Range("F19").Value = 2451
$ws1.Range("F21").Value = 217
$ws1.Range("F23").Value = 454
$ws1.Range("F24").Value = 520
$ws1.Range("F28").Value = 64
$ws1.Range("F29").Value = 1629
$ws1.Range("F30").Value = 241
$ws1.Range("F31").Value = 1143
$ws1.Range("F32").Value = 1911
$ws1.Range("F35").Value = 160
$ws1.Range("F36").Value = 264
$ws1.Range("F39").Value = 313
$ws1.Range("F41").Value = 208

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3088
$ws4.Range("F8").Value = 1687
$ws4.Range("F9").Value = 2102
$ws4.Range("F16").Value = 443
$ws4.Range("F17").Value = 1147
$ws4.Range("F21").Value = 7509
$ws4.Range("F22").Value = 317
$ws4.Range("F23").Value = 2451
$ws4.Range("F26").Value = 217
$ws4.Range("F28").Value = 454
$ws4.Range("F29").Value = 520
$ws4.Range("F33").Value = 64
$ws4.Range("F34").Value = 1629
$ws4.Range("F35").Value = 241
$ws4.Range("F36").Value = 1143
$ws4.Range("F37").Value = 1911
$ws4.Range("F40").Value = 160
$ws4.Range("F41").Value = 264
$ws4.Range("F44").Value = 313
$ws4.Range("F49").Value = 208
